$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 (shifts existing rows 80..107 down to 81..108)
$ws.Rows("80").Insert()

# Populate the newly inserted row 80 with the new weekly price record
$ws.Range("A80").Value = 7
$ws.Range("B80").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C80").Value = "Ñuble"
$ws.Range("D80").Value = 44524
$ws.Range("E80").Value = 16
$ws.Range("F80").Value = 100112045
$ws.Range("G80").Value = "Zapallo"
$ws.Range("H80").Value = "Paine"
$ws.Range("I80").Value = "1a (guarda)"
$ws.Range("J80").Value = 160
$ws.Range("K80").Value = 220
$ws.Range("L80").Value = 250
$ws.Range("M80").Value = 235
$ws.Range("N80").Value = "$/kilo (volumen en unidades)"
$ws.Range("O80").Value = "Región de O'Higgins"
$ws.Range("P80").Value = 235
$ws.Range("Q80").Value = 1
$ws.Range("R80").Value = "Hortaliza"
